# Excel COM-interop edit script
# 1) Table1_Summary_Statistics!B2: "Studies analyzed" value 51 -> 49 (kept as text)
# 2) Raw_Summary_Statistics: insert a new first column "N_observations" (=49)
#    shifting all existing columns one to the right (N_studies..Median_total_variables)

$wb = $excel.ActiveWorkbook

# --- 1) Table1_Summary_Statistics ---------------------------------------
# B2 ("Studies analyzed") holds a number-looking label as TEXT, not a
# number; force text formatting before writing so "49" isn't
# auto-converted to the numeric value 49.
$ws1 = $wb.Worksheets.Item("Table1_Summary_Statistics")
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "49"

# --- 2) Raw_Summary_Statistics ------------------------------------------
$ws9 = $wb.Worksheets.Item("Raw_Summary_Statistics")

# Insert a new column before column A; existing data (incl. formatting)
# shifts right from A:S to B:T.
$ws9.Columns.Item(1).Insert()

# Give the new header cell (A1) the same look as the rest of the header row.
$ws9.Range("B1").Copy()
$ws9.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws9.Range("A1").Value = "N_observations"
$ws9.Range("A2").Value = 49

# Re-assert the shifted numeric cells on row 2 with their exact literal
# values so the round trip through the column-insert operation doesn't
# introduce floating point representation drift.
$ws9.Range("B2").Value = 51
$ws9.Range("C2").Value = 9
$ws9.Range("D2").Value = 27
$ws9.Range("E2").Value = 1.63293
$ws9.Range("F2").Value = 1.2
$ws9.Range("G2").Value = 1.91072
$ws9.Range("H2").Value = 0.000136
$ws9.Range("I2").Value = 8.48
$ws9.Range("J2").Value = 0.265
$ws9.Range("K2").Value = 2.63
$ws9.Range("L2").Value = 2.365
$ws9.Range("M2").Value = 2.108
$ws9.Range("N2").Value = 7.929
$ws9.Range("O2").Value = 2003
$ws9.Range("P2").Value = 2025
$ws9.Range("Q2").Value = 22
$ws9.Range("R2").Value = 94.09999999999999
$ws9.Range("S2").Value = 21.9
$ws9.Range("T2").Value = 21
